# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.707.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '''1.849.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  -2.32%  '
$ws.Range("D5").Value = '''320.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("E6").Value = '  -2.14%  '
$ws.Range("D7").Value = '''0.4309'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.52%  '
$ws.Range("D8").Value = '''0.3741'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("D9").Value = '''0.07359'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.32%  '
$ws.Range("D10").Value = '''0.8806'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").Value = '''21.65'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").Value = '''1.862.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("D13").Value = '''6.735'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = '''5.454'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").Value = '''0.07155'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("D16").Value = '''87.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.64%  '
$ws.Range("E17").Value = '  -2.40%  '
$ws.Range("D18").Value = '''0.000009002'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").Value = '''15.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("D21").Value = '''27.720.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = '''5.250'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("D23").Value = '''11.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").Value = '''2.085.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("D25").Value = '''2.006'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("D26").Value = '''155.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.80%  '
$ws.Range("D27").Value = '''18.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("D28").Value = '''2.131'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.94%  '
$ws.Range("D29").Value = '''5.384'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.16%  '
$ws.Range("D30").Value = '''119.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("D31").Value = '''0.08956'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("D32").Value = '''1.240'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("D33").Value = '''0.7783'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("D34").Value = '''4.570'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").Value = '''2.925'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.91%  '
$ws.Range("E36").Value = '  -2.27%  '
$ws.Range("D37").Value = '''1.138'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("D38").Value = '''0.05344'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = '''0.01974'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("D40").Value = '''7.261'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.54%  '
$ws.Range("D41").Value = '''2.871'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("D42").Value = '''0.5156'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = '''0.1684'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").Value = '''8.812'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.50%  '
$ws.Range("D45").Value = '''109.55'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").Value = '''10.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").Value = '''0.4734'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''1.703'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.10%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.06496'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("E50").Value = '  -2.31%  '
$ws.Range("D51").Value = '''1.874'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.94%  '
